$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.022.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.788.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("E6").Value = "  +0.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5205"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3795"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07828"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.84%  "

$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.088"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.22%  "

$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("C12").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.006"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.247"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.794.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.251"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001078"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.55%  "

$ws.Range("E19").Value = "  -2.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.930"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.054.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.05%  "

$ws.Range("E25").Value = "  -0.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.993.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.309"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1067"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.042"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.677"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.501"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07192"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02309"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.735"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2124"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.034"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6103"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.161"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.365"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.756"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5896"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.221"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.907"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06727"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.33%  "

Write-Host "Updated cryptos list"